$d = $word.ActiveDocument

# Widen the third (results) column of the performance table to match the
# updated, slightly longer confidence-interval strings.
$d.Tables(1).Columns(3).Width = 150.75

# Refresh the imputed-performance figures (2-year horizon).
$d.Content.Find.Execute("1.24%", $false, $false, $false, $false, $false, $true, 1, $false, "1.94%", 2) | Out-Null
$d.Content.Find.Execute("2.6% (2.41% to 2.78%)", $false, $false, $false, $false, $false, $true, 1, $false, "2.73% (2.54% to 2.92%)", 2) | Out-Null
$d.Content.Find.Execute("2.09 (1.92 to 2.27)", $false, $false, $false, $false, $false, $true, 1, $false, "1.41 (1.29 to 1.54)", 2) | Out-Null
$d.Content.Find.Execute("1.35% (1.16% to 1.55%)", $false, $false, $false, $false, $false, $true, 1, $false, "0.79% (0.58% to 1.01%)", 2) | Out-Null
$d.Content.Find.Execute("0.63 (0.51 to 0.75)", $false, $false, $false, $false, $false, $true, 1, $false, "-0.08 (-0.22 to 0.07)", 2) | Out-Null
$d.Content.Find.Execute("-0.26 (-0.33 to -0.19)", $false, $false, $false, $false, $false, $true, 1, $false, "-0.37 (-0.46 to -0.29)", 2) | Out-Null
$d.Content.Find.Execute("0.91 (0.9 to 0.92)", $false, $false, $false, $false, $false, $true, 1, $false, "0.9 (0.89 to 0.92)", 2) | Out-Null
$d.Content.Find.Execute("0.02 (0.02 to 0.02)", $false, $false, $false, $false, $false, $true, 1, $false, "0.02 (0.02 to 0.03)", 2) | Out-Null

# Refresh the imputed-performance figures (5-year horizon).
$d.Content.Find.Execute("3.94%", $false, $false, $false, $false, $false, $true, 1, $false, "5.4%", 2) | Out-Null
$d.Content.Find.Execute("4.49% (4.24% to 4.74%)", $false, $false, $false, $false, $false, $true, 1, $false, "4.76% (4.51% to 5.02%)", 2) | Out-Null
$d.Content.Find.Execute("1.14 (1.06 to 1.22)", $false, $false, $false, $false, $false, $true, 1, $false, "0.88 (0.82 to 0.94)", 2) | Out-Null
$d.Content.Find.Execute("0.54% (0.24% to 0.85%)", $false, $false, $false, $false, $false, $true, 1, $false, "-0.64% (-0.98% to -0.29%)", 2) | Out-Null
$d.Content.Find.Execute("-0.07 (-0.17 to 0.04)", $false, $false, $false, $false, $false, $true, 1, $false, "-0.59 (-0.68 to -0.5)", 2) | Out-Null
$d.Content.Find.Execute("-0.28 (-0.33 to -0.23)", $false, $false, $false, $false, $false, $true, 1, $false, "-0.38 (-0.43 to -0.33)", 2) | Out-Null
$d.Content.Find.Execute("0.89 (0.88 to 0.9)", $false, $false, $false, $false, $false, $true, 1, $false, "0.88 (0.87 to 0.89)", 2) | Out-Null
$d.Content.Find.Execute("0.04 (0.03 to 0.04)", $false, $false, $false, $false, $false, $true, 1, $false, "0.04 (0.04 to 0.04)", 2) | Out-Null
